# Insert a new pricing record for "Macroferia Regional de Talca - Zanahoria"
# as a new row 421, shifting the existing rows 421:447 down to 422:448.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("421:421").Insert()

$ws.Cells.Item(421, 1).Value = 5
$ws.Cells.Item(421, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(421, 3).Value = "Maule"
$ws.Cells.Item(421, 4).Value = 44931
$ws.Cells.Item(421, 5).Value = 7
$ws.Cells.Item(421, 6).Value = 100114013
$ws.Cells.Item(421, 7).Value = "Zanahoria"
$ws.Cells.Item(421, 8).Value = "Sin especificar"
$ws.Cells.Item(421, 9).Value = "Primera"
$ws.Cells.Item(421, 10).Value = 400
$ws.Cells.Item(421, 11).Value = 12000
$ws.Cells.Item(421, 12).Value = 12000
$ws.Cells.Item(421, 13).Value = 12000
$ws.Cells.Item(421, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(421, 15).Value = "Región de Ñuble"
$ws.Cells.Item(421, 16).Value = 600
$ws.Cells.Item(421, 17).Value = 20
$ws.Cells.Item(421, 18).Value = "Hortaliza"
